# Amend the input distribution parameter strings on each patient-trajectory
# sheet: every "lognorm,<mean>" string becomes "lognorm,<mean>,<sd>" (sd is
# 10% of the mean), and a couple of stale "lognorm,6"/"lognorm,4.8,..." typos
# are corrected to the canonical "lognorm,0.6,0.06" value used elsewhere.
#
# Layout of each "trajN" sheet is Area (col A) / distribution (col B) pairs,
# starting at row 2 (row 1 is the header "Area" / "groupA").

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("traj1")
$ws2 = $wb.Worksheets.Item("traj2")
$ws3 = $wb.Worksheets.Item("traj3")
$ws4 = $wb.Worksheets.Item("traj4")
$ws5 = $wb.Worksheets.Item("traj5")

# traj1: ECU  lognorm,6   -> lognorm,0.6,0.06
$ws1.Range("B2").Value = "lognorm,0.6,0.06"

# traj2: ICU  lognorm,6   -> lognorm,0.6,0.06
$ws2.Range("B2").Value = "lognorm,0.6,0.06"

# traj3: ECU  lognorm,0.6 -> lognorm,0.6,0.06
#        ICU  lognorm,5.4 -> lognorm,5.4,0.54
$ws3.Range("B2").Value = "lognorm,0.6,0.06"
$ws3.Range("B3").Value = "lognorm,5.4,0.54"

# traj4: ICU  lognorm,5.4 -> lognorm,5.4,0.54
#        ECU  lognorm,0.6 -> lognorm,0.6,0.06
$ws4.Range("B2").Value = "lognorm,5.4,0.54"
$ws4.Range("B3").Value = "lognorm,0.6,0.06"

# traj5: ECU  lognorm,0.6 -> lognorm,0.6,0.06
#        ICU  lognorm,4.8 -> lognorm,4.8,0.48
#        ECU  lognorm,0.6 -> lognorm,0.6,0.06
$ws5.Range("B2").Value = "lognorm,0.6,0.06"
$ws5.Range("B3").Value = "lognorm,4.8,0.48"
$ws5.Range("B4").Value = "lognorm,0.6,0.06"

# Update each sheet's stored selection / active cell to match the edited
# workbook, then restore traj1 as the active tab (it was active before the
# edit and must remain so).
$ws1.Range("B2").Select()
$ws2.Range("B2").Select()
$ws3.Range("B2").Select()
$ws4.Range("B3").Select()
$ws5.Range("B4").Select()

$ws1.Activate()
$ws1.Range("B2").Select()
